$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Cells.Item(5, 7).Value = 2.63  # G5: 2.55 -> 2.63
$ws.Cells.Item(5, 8).Value = 2.7  # H5: 2.8 -> 2.7
$ws.Cells.Item(5, 10).Value = 1.18  # J5: 1.17 -> 1.18
$ws.Cells.Item(5, 11).Value = 4.5  # K5: 5 -> 4.5
$ws.Cells.Item(5, 12).Value = 1.83  # L5: 1.73 -> 1.83
$ws.Cells.Item(5, 13).Value = 1.83  # M5: 2 -> 1.83
$ws.Cells.Item(5, 14).Value = 3.6  # N5: 3.5 -> 3.6
$ws.Cells.Item(5, 15).Value = 1.29  # O5: 1.3 -> 1.29
$ws.Cells.Item(5, 16).Value = 1.83  # P5: 1.78 -> 1.83
$ws.Cells.Item(5, 17).Value = 1.98  # Q5: 2.03 -> 1.98
$ws.Cells.Item(5, 18).Value = 2.75  # R5: 2.63 -> 2.75
$ws.Cells.Item(5, 19).Value = 1.4  # S5: 1.44 -> 1.4
$ws.Cells.Item(5, 20).Value = 5  # T5: 5.5 -> 5
$ws.Cells.Item(5, 21).Value = 11  # U5: 10 -> 11
$ws.Cells.Item(5, 22).Value = 13  # V5: 12 -> 13
$ws.Cells.Item(5, 23).Value = 29  # W5: 26 -> 29
$ws.Cells.Item(5, 24).Value = 34  # X5: 29 -> 34
$ws.Cells.Item(5, 26).Value = 4.5  # Z5: 5 -> 4.5
$ws.Cells.Item(5, 28).Value = 26  # AB5: 23 -> 26
$ws.Cells.Item(5, 29).Value = 126  # AC5: 101 -> 126
$ws.Cells.Item(5, 31).Value = 5.5  # AE5: 6 -> 5.5
$ws.Cells.Item(5, 33).Value = 15  # AG5: 13 -> 15
$ws.Cells.Item(5, 34).Value = 34  # AH5: 41 -> 34

# Row 8
$ws.Cells.Item(8, 8).Value = 3.4  # H8: 3.45 -> 3.4
$ws.Cells.Item(8, 9).Value = 3.9  # I8: 3.8 -> 3.9
$ws.Cells.Item(8, 11).Value = 7.8  # K8: 7.9 -> 7.8
$ws.Cells.Item(8, 12).Value = 1.29  # L8: 1.28 -> 1.29
$ws.Cells.Item(8, 13).Value = 3.4  # M8: 3.45 -> 3.4
$ws.Cells.Item(8, 14).Value = 1.87  # N8: 1.85 -> 1.87
$ws.Cells.Item(8, 15).Value = 1.88  # O8: 1.9 -> 1.88
$ws.Cells.Item(8, 16).Value = 1.42  # P8: 1.4 -> 1.42
$ws.Cells.Item(8, 17).Value = 2.75  # Q8: 2.8 -> 2.75
$ws.Cells.Item(8, 20).Value = 7.5  # T8: 7.6 -> 7.5
$ws.Cells.Item(8, 23).Value = 19  # W8: 18.5 -> 19
$ws.Cells.Item(8, 26).Value = 7.8  # Z8: 7.9 -> 7.8
$ws.Cells.Item(8, 27).Value = 6.8  # AA8: 7 -> 6.8
$ws.Cells.Item(8, 28).Value = 14.5  # AB8: 15 -> 14.5
$ws.Cells.Item(8, 32).Value = 23  # AF8: 22 -> 23
$ws.Cells.Item(8, 34).Value = 65  # AH8: 60 -> 65
$ws.Cells.Item(8, 35).Value = 40  # AI8: 37 -> 40

# Row 9
$ws.Cells.Item(9, 7).Value = 2  # G9: 2.22 -> 2
$ws.Cells.Item(9, 8).Value = 3.2  # H9: 3.15 -> 3.2
$ws.Cells.Item(9, 9).Value = 3.9  # I9: 3.3 -> 3.9
$ws.Cells.Item(9, 10).Value = 1.09  # J9: 1.1 -> 1.09
$ws.Cells.Item(9, 11).Value = 6.4  # K9: 6.3 -> 6.4
$ws.Cells.Item(9, 12).Value = 1.42  # L9: 1.44 -> 1.42
$ws.Cells.Item(9, 13).Value = 2.7  # M9: 2.65 -> 2.7
$ws.Cells.Item(9, 14).Value = 2.25  # N9: 2.3 -> 2.25
$ws.Cells.Item(9, 15).Value = 1.6  # O9: 1.57 -> 1.6
$ws.Cells.Item(9, 20).Value = 6  # T9: 6.4 -> 6
$ws.Cells.Item(9, 21).Value = 9.25  # U9: 10.5 -> 9.25
$ws.Cells.Item(9, 22).Value = 9.25  # V9: 10 -> 9.25
$ws.Cells.Item(9, 23).Value = 19  # W9: 23 -> 19
$ws.Cells.Item(9, 24).Value = 19.5  # X9: 22 -> 19.5
$ws.Cells.Item(9, 26).Value = 6.4  # Z9: 6.3 -> 6.4
$ws.Cells.Item(9, 27).Value = 6.5  # AA9: 6.4 -> 6.5
$ws.Cells.Item(9, 31).Value = 9  # AE9: 7.9 -> 9
$ws.Cells.Item(9, 32).Value = 21  # AF9: 17 -> 21
$ws.Cells.Item(9, 33).Value = 14.5  # AG9: 13 -> 14.5
$ws.Cells.Item(9, 34).Value = 70  # AH9: 50 -> 70
$ws.Cells.Item(9, 35).Value = 45  # AI9: 37 -> 45
$ws.Cells.Item(9, 36).Value = 60  # AJ9: 50 -> 60

# Row 10
$ws.Cells.Item(10, 14).Value = 2.4  # N10: 2.5 -> 2.4
$ws.Cells.Item(10, 15).Value = 1.53  # O10: 1.5 -> 1.53
$ws.Cells.Item(10, 16).Value = 1.53  # P10: 1.57 -> 1.53
$ws.Cells.Item(10, 17).Value = 2.38  # Q10: 2.25 -> 2.38
$ws.Cells.Item(10, 20).Value = 6.5  # T10: 6 -> 6.5
$ws.Cells.Item(10, 24).Value = 21  # X10: 23 -> 21
$ws.Cells.Item(10, 26).Value = 7  # Z10: 6.5 -> 7
$ws.Cells.Item(10, 31).Value = 8  # AE10: 7.5 -> 8
$ws.Cells.Item(10, 33).Value = 12  # AG10: 13 -> 12
$ws.Cells.Item(10, 35).Value = 29  # AI10: 34 -> 29

# Row 40
$ws.Cells.Item(40, 7).Value = 1.6  # G40: 1.73 -> 1.6
$ws.Cells.Item(40, 8).Value = 3.4  # H40: 3.25 -> 3.4
$ws.Cells.Item(40, 9).Value = 6  # I40: 5.25 -> 6
$ws.Cells.Item(40, 12).Value = 1.33  # L40: 1.36 -> 1.33
$ws.Cells.Item(40, 13).Value = 3.25  # M40: 3 -> 3.25
$ws.Cells.Item(40, 14).Value = 2.1  # N40: 2.15 -> 2.1
$ws.Cells.Item(40, 15).Value = 1.7  # O40: 1.67 -> 1.7
$ws.Cells.Item(40, 18).Value = 2.1  # R40: 2 -> 2.1
$ws.Cells.Item(40, 19).Value = 1.67  # S40: 1.73 -> 1.67
$ws.Cells.Item(40, 21).Value = 7  # U40: 7.5 -> 7
$ws.Cells.Item(40, 23).Value = 12  # W40: 13 -> 12
$ws.Cells.Item(40, 27).Value = 7  # AA40: 6.5 -> 7
$ws.Cells.Item(40, 28).Value = 19  # AB40: 17 -> 19
$ws.Cells.Item(40, 31).Value = 13  # AE40: 12 -> 13
$ws.Cells.Item(40, 32).Value = 29  # AF40: 26 -> 29
$ws.Cells.Item(40, 33).Value = 19  # AG40: 17 -> 19
$ws.Cells.Item(40, 34).Value = 67  # AH40: 51 -> 67
$ws.Cells.Item(40, 35).Value = 51  # AI40: 41 -> 51
